$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: Taxonsorteringsordning bumped 79243 -> 79244 ---
$ws.Range("B3").Value = 79244

# --- Rows 4 & 5: the "Antal"/"Enhet" observation (I/J) that used to sit on
#     row 5 moves to row 4 (row 4 previously had no Antal/Enhet at all).
#     Use Copy so the original text-cell typing (e.g. "10") is preserved
#     instead of Excel auto-coercing a numeric-looking string to a number. ---
$ws.Range("I5").Copy($ws.Range("I4"))
$ws.Range("J5").Copy($ws.Range("J4"))
$ws.Range("I5").ClearContents()
$ws.Range("J5").ClearContents()

# --- Rows 4 & 5: the Id values (A) swap between the two rows ---
$a4 = $ws.Range("A4").Value()
$a5 = $ws.Range("A5").Value()
$ws.Range("A4").Value = $a5
$ws.Range("A5").Value = $a4

# --- Rows 4 & 5: Taxonsorteringsordning bumped 79243 -> 79244 on both ---
$ws.Range("B4").Value = 79244
$ws.Range("B5").Value = 79244

# --- Rows 4 & 5: coordinates (Ost/Nord) swap between the two rows ---
$q4 = $ws.Range("Q4").Value()
$r4 = $ws.Range("R4").Value()
$q5 = $ws.Range("Q5").Value()
$r5 = $ws.Range("R5").Value()
$ws.Range("Q4").Value = $q5
$ws.Range("R4").Value = $r5
$ws.Range("Q5").Value = $q4
$ws.Range("R5").Value = $r4

# --- Rows 4 & 5: Externid / Starttid / Sluttid swap between the two rows ---
$x4 = $ws.Range("X4").Value()
$z4 = $ws.Range("Z4").Value()
$ab4 = $ws.Range("AB4").Value()
$x5 = $ws.Range("X5").Value()
$z5 = $ws.Range("Z5").Value()
$ab5 = $ws.Range("AB5").Value()
$ws.Range("X4").Value = $x5
$ws.Range("Z4").Value = $z5
$ws.Range("AB4").Value = $ab5
$ws.Range("X5").Value = $x4
$ws.Range("Z5").Value = $z4
$ws.Range("AB5").Value = $ab4

# --- Row 6: Taxonsorteringsordning bumped 79862 -> 79863 ---
$ws.Range("B6").Value = 79863
